$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the style of an existing header cell (H1) to the new
# header cells first, then set their values/text.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF)
$values = @(
    @(5, 6),
    @(9, 9),
    @(9, 9),
    @(6, 7),
    @(8, 9),
    @(9, 9),
    @(6, 7),
    @(5, 5),
    @(8, 8),
    @(6, 6),
    @(3, 3),
    @(9, 9),
    @(5, 6),
    @(9, 9),
    @(9, 9),
    @(9, 9),
    @(4, 4)
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
